# Append a new scraped listing row (2025-11-01 01:23:01) to the
# "ランサーズ" sheet and refresh the "取得日時" timestamp on every
# existing row to the same collection time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-01 01:23:01"

# Update the acquisition timestamp (column A) for the existing rows 2-9.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Append the new row of data at row 10.
$row = 10
$ws.Cells.Item($row, 1).Value = $newTimestamp
$ws.Cells.Item($row, 2).Value = "【電子工作】基盤にDCケーブルのターミナルと抵抗を追加したい方募集!"
$ws.Cells.Item($row, 3).Value = "システム開発"
$ws.Cells.Item($row, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item($row, 5).Value = "期限情報なし"
$ws.Cells.Item($row, 6).Value = "https://www.lancers.jp/work/detail/5424906"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 6), "https://www.lancers.jp/work/detail/5424906")
$ws.Cells.Item($row, 6).Style = "Hyperlink"
$ws.Cells.Item($row, 7).Value = 10
